# Weekly update: insert a new price record for "Vega Modelo de Temuco - Zanahoria"
# at row 393, pushing the existing rows 393:456 down to 394:457.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 393 (shifts 393:456 -> 394:457, and extends the
# used range / dimension to A1:R457 automatically, inheriting the date-column
# number format from the row above for the new D393 cell).
$ws.Rows.Item(393).Insert()

# Populate the newly inserted row with this week's new record.
$ws.Range("A393").Value = 10
$ws.Range("B393").Value = "Vega Modelo de Temuco"
$ws.Range("C393").Value = "La Araucanía"
$ws.Range("D393").Value = 45034
$ws.Range("E393").Value = 9
$ws.Range("F393").Value = 100114013
$ws.Range("G393").Value = "Zanahoria"
$ws.Range("H393").Value = "Sin especificar"
$ws.Range("I393").Value = "Primera"
$ws.Range("J393").Value = 110
$ws.Range("K393").Value = 6000
$ws.Range("L393").Value = 6000
$ws.Range("M393").Value = 6000
$ws.Range("N393").Value = "`$/saco 25 kilos"
$ws.Range("O393").Value = "Región de La Araucanía"
$ws.Range("P393").Value = 240
$ws.Range("Q393").Value = 25
$ws.Range("R393").Value = "Hortaliza"
